$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3799.5
$ws.Range("H4").Value = 608.3333
$ws.Range("I4").Value = 600
$ws.Range("K4").Value = 600
$ws.Range("M4").Value = -486
$ws.Range("H12").Value = 1827.0454
$ws.Range("I12").Value = 1846.8235
$ws.Range("K12").Value = 1846.8235
$ws.Range("M12").Value = -1676.8235
$ws.Range("H32").Value = 5463.8
$ws.Range("I32").Value = 4148
$ws.Range("K32").Value = 4148
$ws.Range("M32").Value = -3822
$ws.Range("H94").Value = 2211.2
$ws.Range("I94").Value = 2211.2
$ws.Range("K94").Value = 2211.2
$ws.Range("M94").Value = -1760.2
$ws.Range("H116").Value = 4564.077
$ws.Range("I116").Value = 4304.875
$ws.Range("J116").Value = 4978.8
$ws.Range("K116").Value = 4304.875
$ws.Range("L116").Value = 4978.8
$ws.Range("M116").Value = -862.875
$ws.Range("N116").Value = -11862.8
$ws.Range("H132").Value = 2098.1228
$ws.Range("I132").Value = 1073.4166
$ws.Range("K132").Value = 3220.2498
$ws.Range("M132").Value = -690.2498000000001
$ws.Range("H137").Value = 2793.5
$ws.Range("I137").Value = 764.7273
$ws.Range("J137").Value = 4822.273
$ws.Range("K137").Value = 2294.1819
$ws.Range("L137").Value = 14466.819
$ws.Range("M137").Value = 255.8181
$ws.Range("N137").Value = -19566.819
$ws.Range("H138").Value = 2809.246
$ws.Range("J138").Value = 3044.75
$ws.Range("L138").Value = 9134.25
$ws.Range("N138").Value = -19414.25
$ws.Range("H140").Value = 21780
$ws.Range("J140").Value = 21780
$ws.Range("L140").Value = 21780
$ws.Range("N140").Value = -32140

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 125570.31
$ws.Range("I4").Value = 167207.75
$ws.Range("K4").Value = 167207.75
$ws.Range("M4").Value = -167091.75
$ws.Range("H5").Value = 701.9231
$ws.Range("I5").Value = 747.2222
$ws.Range("K5").Value = 747.2222
$ws.Range("M5").Value = -635.2222
$ws.Range("H32").Value = 5688.283
$ws.Range("I32").Value = 4511.074
$ws.Range("K32").Value = 4511.074
$ws.Range("M32").Value = -4224.074
$ws.Range("H45").Value = 10310.533
$ws.Range("I45").Value = 15478.75
$ws.Range("K45").Value = 15478.75
$ws.Range("M45").Value = -15101.75
$ws.Range("H61").Value = 1397.6364
$ws.Range("I61").Value = 1295.4117
$ws.Range("K61").Value = 1295.4117
$ws.Range("M61").Value = -1083.4117
$ws.Range("H74").Value = 1932.9546
$ws.Range("I74").Value = 1122
$ws.Range("J74").Value = 3104.3333
$ws.Range("K74").Value = 1122
$ws.Range("L74").Value = 3104.3333
$ws.Range("M74").Value = -248
$ws.Range("N74").Value = -4852.3333
$ws.Range("H77").Value = 1932.9546
$ws.Range("I77").Value = 1122
$ws.Range("J77").Value = 3104.3333
$ws.Range("K77").Value = 5610
$ws.Range("L77").Value = 15521.6665
$ws.Range("M77").Value = -1242
$ws.Range("N77").Value = -24257.6665
$ws.Range("H132").Value = 1954.569
$ws.Range("I132").Value = 1766.6666
$ws.Range("J132").Value = 2605
$ws.Range("K132").Value = 5299.9998
$ws.Range("L132").Value = 7815
$ws.Range("M132").Value = -2769.9998
$ws.Range("N132").Value = -12875
$ws.Range("H136").Value = 1397.6364
$ws.Range("I136").Value = 1295.4117
$ws.Range("K136").Value = 3886.2351
$ws.Range("M136").Value = -1336.2351

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 701.9231
$ws.Range("I4").Value = 747.2222
$ws.Range("K4").Value = 747.2222
$ws.Range("M4").Value = -632.2222
$ws.Range("H43").Value = 130000
$ws.Range("J43").Value = 130000
$ws.Range("L43").Value = 130000
$ws.Range("N43").Value = -130362
$ws.Range("H86").Value = 2031.6923
$ws.Range("I86").Value = 1927.7273
$ws.Range("K86").Value = 1927.7273
$ws.Range("M86").Value = -804.7273
$ws.Range("H89").Value = 2031.6923
$ws.Range("I89").Value = 1927.7273
$ws.Range("K89").Value = 9638.636500000001
$ws.Range("M89").Value = -4022.636500000001
$ws.Range("H94").Value = 66669076
$ws.Range("I94").Value = 111113020
$ws.Range("J94").Value = 3153.3333
$ws.Range("K94").Value = 111113020
$ws.Range("L94").Value = 3153.3333
$ws.Range("M94").Value = -111112569
$ws.Range("N94").Value = -4055.3333
$ws.Range("H105").Value = 2137.7693
$ws.Range("I105").Value = 2132.625
$ws.Range("K105").Value = 2132.625
$ws.Range("M105").Value = -385.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66666956
$ws.Range("I7").Value = 90909370
$ws.Range("K7").Value = 90909370
$ws.Range("M7").Value = -90909257
$ws.Range("H58").Value = 2479.9443
$ws.Range("I58").Value = 2282.1428
$ws.Range("J58").Value = 2605.818
$ws.Range("K58").Value = 2282.1428
$ws.Range("L58").Value = 2605.818
$ws.Range("M58").Value = -2079.1428
$ws.Range("N58").Value = -3011.818
$ws.Range("H136").Value = 2479.9443
$ws.Range("I136").Value = 2282.1428
$ws.Range("J136").Value = 2605.818
$ws.Range("K136").Value = 6846.428400000001
$ws.Range("L136").Value = 7817.454000000001
$ws.Range("M136").Value = -4296.428400000001
$ws.Range("N136").Value = -12917.454

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 6121.647
$ws.Range("I106").Value = 3999
$ws.Range("J106").Value = 6404.6665
$ws.Range("K106").Value = 11997
$ws.Range("L106").Value = 19213.9995
$ws.Range("M106").Value = -11051
$ws.Range("N106").Value = -21105.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6456
$ws.Range("I70").Value = 6657.4165
$ws.Range("J70").Value = 6153.875
$ws.Range("K70").Value = 6657.4165
$ws.Range("L70").Value = 6153.875
$ws.Range("M70").Value = -6387.4165
$ws.Range("N70").Value = -6693.875
$ws.Range("H73").Value = 6456
$ws.Range("I73").Value = 6657.4165
$ws.Range("J73").Value = 6153.875
$ws.Range("K73").Value = 6657.4165
$ws.Range("L73").Value = 6153.875
$ws.Range("M73").Value = -5721.4165
$ws.Range("N73").Value = -8025.875
$ws.Range("H132").Value = 4189.615
$ws.Range("I132").Value = 2834.0322
$ws.Range("J132").Value = 9442.5
$ws.Range("K132").Value = 8502.096600000001
$ws.Range("L132").Value = 28327.5
$ws.Range("M132").Value = -5972.096600000001
$ws.Range("N132").Value = -33387.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 17393.445
$ws.Range("I40").Value = 8757.111000000001
$ws.Range("K40").Value = 8757.111000000001
$ws.Range("M40").Value = -8621.111000000001
$ws.Range("H68").Value = 4390.467
$ws.Range("I68").Value = 4280.4614
$ws.Range("K68").Value = 4280.4614
$ws.Range("M68").Value = -3531.4614
$ws.Range("H71").Value = 4390.467
$ws.Range("I71").Value = 4280.4614
$ws.Range("K71").Value = 21402.307
$ws.Range("M71").Value = -17658.307
$ws.Range("H82").Value = 4462.4614
$ws.Range("I82").Value = 5138.25
$ws.Range("J82").Value = 4162.1113
$ws.Range("K82").Value = 5138.25
$ws.Range("L82").Value = 4162.1113
$ws.Range("M82").Value = -4777.25
$ws.Range("N82").Value = -4884.1113
$ws.Range("H85").Value = 4462.4614
$ws.Range("I85").Value = 5138.25
$ws.Range("J85").Value = 4162.1113
$ws.Range("K85").Value = 5138.25
$ws.Range("L85").Value = 4162.1113
$ws.Range("M85").Value = -3890.25
$ws.Range("N85").Value = -6658.1113
$ws.Range("H136").Value = 7017.2905
$ws.Range("I136").Value = 8943.532999999999
$ws.Range("K136").Value = 26830.599
$ws.Range("M136").Value = -24280.599

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1164889.5
$ws.Range("I132").Value = 1521544.1
$ws.Range("J132").Value = 5762.125
$ws.Range("K132").Value = 4564632.300000001
$ws.Range("L132").Value = 17286.375
$ws.Range("M132").Value = -4562102.300000001
$ws.Range("N132").Value = -22346.375
$ws.Range("H136").Value = 1133.3846
$ws.Range("I136").Value = 1024.375
$ws.Range("J136").Value = 1307.8
$ws.Range("K136").Value = 3073.125
$ws.Range("L136").Value = 3923.4
$ws.Range("M136").Value = -523.125
$ws.Range("N136").Value = -9023.4
